# Applies the "Add files via upload" edit to the "classifica" block at the
# bottom of the single worksheet:
#   - tweaks a handful of score values in the existing ranking table
#     (rows 100:103) and freezes the SUM() formulas in column L into
#     plain numbers
#   - appends a small "tot UBUNTU" / "tot WINDOWS" summary table
#     (rows 104:108)
#   - leaves a formatted-but-empty cell at G116 (underlined font +
#     the workbook's custom decimal number format) and selects it,
#     matching where the author's cursor ended up

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- corrected scores inside the existing ranking table -------------
$ws.Range("K100").Value = 3
$ws.Range("E101").Value = 2
$ws.Range("G101").Value = 3
$ws.Range("G102").Value = 2
$ws.Range("K102").Value = 2
$ws.Range("M102").Value = 2
$ws.Range("E103").Value = 1
$ws.Range("G103").Value = 1
$ws.Range("K103").Value = 1

# ---- column L held =SUM(D:K) formulas; pin down the refreshed totals
$ws.Range("L100").Value = 31
$ws.Range("L101").Value = 24
$ws.Range("L102").Value = 15
$ws.Range("L103").Value = 10

# ---- new "tot UBUNTU" / "tot WINDOWS" mini summary table -------------
$ws.Range("G104").Value = "tot UBUNTU"
$ws.Range("K104").Value = "tot WINDOWS"

$ws.Range("F105").Value = "c++"
$ws.Range("G105").Value = 16
$ws.Range("H105").Value = 4
$ws.Range("K105").Value = 15
$ws.Range("L105").Value = 4

$ws.Range("F106").Value = "python"
$ws.Range("G106").Value = 11
$ws.Range("H106").Value = 3
$ws.Range("K106").Value = 13
$ws.Range("L106").Value = 3

$ws.Range("F107").Value = "R"
$ws.Range("G107").Value = 8
$ws.Range("H107").Value = 2
$ws.Range("K107").Value = 7
$ws.Range("L107").Value = 2

$ws.Range("F108").Value = "MATLAB"
$ws.Range("G108").Value = 5
$ws.Range("H108").Value = 1
$ws.Range("K108").Value = 5
$ws.Range("L108").Value = 1

# ---- stray formatted cell further down the sheet ---------------------
$g116 = $ws.Range("G116")
$g116.NumberFormat = "0.0000000000000000"
$g116.Font.Underline = $true

# ---- cursor/selection ends up on the new cell -------------------------
$g116.Select()
